$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Max Weight constraint for GLD (row 8) to 100%
$ws.Range("C8").Value = 1

# Update the active cell selection to reflect where the user clicked last
[void]$ws.Range("J12").Select()
